$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 481.57895
$ws.Range("I19").Value = 575.2857
$ws.Range("J19").Value = 426.91666
$ws.Range("K19").Value = 575.2857
$ws.Range("L19").Value = 426.91666
$ws.Range("M19").Value = -400.2857
$ws.Range("N19").Value = -776.91666

$ws.Range("H28").Value = 892.41174
$ws.Range("I28").Value = 1078.75
$ws.Range("J28").Value = 445.2
$ws.Range("K28").Value = 1078.75
$ws.Range("L28").Value = 445.2
$ws.Range("M28").Value = -593.75
$ws.Range("N28").Value = -1415.2

$ws.Range("H32").Value = 125001200
$ws.Range("I32").Value = 500001000
$ws.Range("J32").Value = 1267
$ws.Range("K32").Value = 500001000
$ws.Range("L32").Value = 1267
$ws.Range("M32").Value = -500000674
$ws.Range("N32").Value = -1919

$ws.Range("H107").Value = 984.08
$ws.Range("I107").Value = 1251.7778
$ws.Range("J107").Value = 295.7143
$ws.Range("K107").Value = 1251.7778
$ws.Range("L107").Value = 295.7143
$ws.Range("M107").Value = 668.2221999999999
$ws.Range("N107").Value = -4135.7143

$ws.Range("H138").Value = 2089.6182
$ws.Range("I138").Value = 1665.8
$ws.Range("J138").Value = 2331.8
$ws.Range("K138").Value = 4997.4
$ws.Range("L138").Value = 6995.400000000001
$ws.Range("M138").Value = 142.6000000000004
$ws.Range("N138").Value = -17275.4


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 42650.5
$ws.Range("J101").Value = 42650.5
$ws.Range("L101").Value = 42650.5
$ws.Range("N101").Value = -49140.5

$ws.Range("H112").Value = 17922.334
$ws.Range("J112").Value = 17922.334
$ws.Range("L112").Value = 17922.334
$ws.Range("N112").Value = -20876.334

$ws.Range("H114").Value = 32049.5
$ws.Range("J114").Value = 32049.5
$ws.Range("L114").Value = 32049.5
$ws.Range("N114").Value = -40727.5

$ws.Range("H132").Value = 2186.652
$ws.Range("I132").Value = 1914.75
$ws.Range("J132").Value = 3999.3333
$ws.Range("K132").Value = 5744.25
$ws.Range("L132").Value = 11997.9999
$ws.Range("M132").Value = -3214.25
$ws.Range("N132").Value = -17057.9999


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H112").Value = 29500
$ws.Range("J112").Value = 29500
$ws.Range("L112").Value = 29500
$ws.Range("N112").Value = -32454

$ws.Range("H122").Value = 32580
$ws.Range("J122").Value = 32580
$ws.Range("L122").Value = 32580
$ws.Range("N122").Value = -42380


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2189.1904
$ws.Range("I31").Value = 1445.12
$ws.Range("K31").Value = 1445.12
$ws.Range("M31").Value = -1150.12

$ws.Range("H34").Value = 2189.1904
$ws.Range("I34").Value = 1445.12
$ws.Range("K34").Value = 1445.12
$ws.Range("M34").Value = -1243.12


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5249.375
$ws.Range("I3").Value = 2290.8333
$ws.Range("J3").Value = 14125
$ws.Range("K3").Value = 6872.499899999999
$ws.Range("L3").Value = 42375
$ws.Range("M3").Value = -6760.499899999999
$ws.Range("N3").Value = -42599

$ws.Range("H68").Value = 1907.6154
$ws.Range("I68").Value = 2800
$ws.Range("J68").Value = 1142.7142
$ws.Range("K68").Value = 8400
$ws.Range("L68").Value = 3428.1426
$ws.Range("M68").Value = -7589
$ws.Range("N68").Value = -5050.142599999999

$ws.Range("H71").Value = 1907.6154
$ws.Range("I71").Value = 2800
$ws.Range("J71").Value = 1142.7142
$ws.Range("K71").Value = 25200
$ws.Range("L71").Value = 10284.4278
$ws.Range("M71").Value = -21144
$ws.Range("N71").Value = -18396.4278

$ws.Range("H99").Value = 2583.647
$ws.Range("I99").Value = 653.6667
$ws.Range("J99").Value = 3636.3635
$ws.Range("K99").Value = 1961.0001
$ws.Range("L99").Value = 10909.0905
$ws.Range("M99").Value = 284.9999
$ws.Range("N99").Value = -15401.0905

$ws.Range("H138").Value = 83335500
$ws.Range("I138").Value = 100001240
$ws.Range("J138").Value = 6800
$ws.Range("K138").Value = 300003720
$ws.Range("L138").Value = 20400
$ws.Range("M138").Value = -299998580
$ws.Range("N138").Value = -30680

$ws.Range("H139").Value = 26394.951
$ws.Range("I139").Value = 1509.2858
$ws.Range("J139").Value = 79994.84
$ws.Range("K139").Value = 4527.857400000001
$ws.Range("L139").Value = 239984.52
$ws.Range("M139").Value = 612.1425999999992
$ws.Range("N139").Value = -250264.52

$ws.Range("H141").Value = 7506
$ws.Range("I141").Value = 2898.889
$ws.Range("J141").Value = 14416.667
$ws.Range("K141").Value = 8696.667000000001
$ws.Range("L141").Value = 43250.001
$ws.Range("M141").Value = -3516.667000000001
$ws.Range("N141").Value = -53610.001


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 28888.666
$ws.Range("J103").Value = 28888.666
$ws.Range("L103").Value = 28888.666
$ws.Range("N103").Value = -31232.666

$ws.Range("H126").Value = 5414.263
$ws.Range("I126").Value = 2419.1
$ws.Range("J126").Value = 8742.223
$ws.Range("K126").Value = 7257.299999999999
$ws.Range("L126").Value = 26226.669
$ws.Range("M126").Value = -4787.299999999999
$ws.Range("N126").Value = -31166.669

$ws.Range("H135").Value = 48851.43
$ws.Range("J135").Value = 48851.43
$ws.Range("L135").Value = 48851.43
$ws.Range("N135").Value = -58991.43


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 43046
$ws.Range("J105").Value = 43046
$ws.Range("L105").Value = 43046
$ws.Range("N105").Value = -50034

$ws.Range("H110").Value = 13900
$ws.Range("J110").Value = 13900
$ws.Range("L110").Value = 13900
$ws.Range("N110").Value = -22080

$ws.Range("H136").Value = 1643.8334
$ws.Range("I136").Value = 1234.6945
$ws.Range("J136").Value = 2871.25
$ws.Range("K136").Value = 3704.0835
$ws.Range("L136").Value = 8613.75
$ws.Range("M136").Value = -1154.0835
$ws.Range("N136").Value = -13713.75


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2397.8262
$ws.Range("I122").Value = 2061.1052
$ws.Range("J122").Value = 3997.25
$ws.Range("K122").Value = 6183.3156
$ws.Range("L122").Value = 11991.75
$ws.Range("M122").Value = -3733.3156
$ws.Range("N122").Value = -16891.75

$ws.Range("H136").Value = 983.9211
$ws.Range("I136").Value = 886.7727
$ws.Range("J136").Value = 1117.5
$ws.Range("K136").Value = 2660.3181
$ws.Range("L136").Value = 3352.5
$ws.Range("M136").Value = -110.3181
$ws.Range("N136").Value = -8452.5

